$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5: rawFrequencyCreationServiceUrl
$ws.Range("A5").Value = "rawFrequencyCreationServiceUrl"
$ws.Hyperlinks.Add($ws.Range("B5"), "http://google.com")
$ws.Range("B5").Style = "Hyperlink"

# Row 6: rawVoltageCreationServiceUrl
$ws.Range("A6").Value = "rawVoltageCreationServiceUrl"
$ws.Hyperlinks.Add($ws.Range("B6"), "http://google.com")
$ws.Range("B6").Style = "Hyperlink"

$ws.Range("A6").Select()
